# Cam Control Design.xlsx edit
# - Insert a new "Tick duration (ns)" column into Tableau1 (Speed sheet),
#   positioned right after "Tick duration (s)" (i.e. becomes column D,
#   everything from the old D ("TOF Period (s)") onward shifts right by one).
# - New column formula: =Tableau1[[#This Row],[Tick duration (s)]]*1000000000
# - Highlight the whole "Prescaler idx = 5" row (row 6) with a yellow fill,
#   marking the chosen prescaler for the new servo driver.
# - Move the active selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Speed")
$lo = $ws.ListObjects.Item("Tableau1")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: remember the current content (formula/value) of the columns that
# need to shift right by one slot: D..I (4..9) -> E..J (5..10).
# ---------------------------------------------------------------------------
$lastRow = $lo.Range.Rows.Count
$origFormula = @{}
$origValue = @{}
for ($c = 4; $c -le 9; $c++) {
    $fcol = @()
    $vcol = @()
    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $c)
        $fcol += $cell.Formula
        $vcol += $cell.Value()
    }
    $origFormula[$c] = $fcol
    $origValue[$c] = $vcol
}

# ---------------------------------------------------------------------------
# Step 2: grow the table by one column on the right (A1:I11 -> A1:J11) so
# there is a tenth column to shift data into.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, 10))) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: shift formatting + content from column c into column c+1, working
# right-to-left so we never overwrite a source column before reading it.
# ---------------------------------------------------------------------------
for ($c = 9; $c -ge 4; $c--) {
    $destCol = $c + 1

    $ws.Cells.Item(1, $c).Resize($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item(1, $destCol).Resize($lastRow, 1).PasteSpecial($xlPasteFormats) | Out-Null

    for ($r = 1; $r -le $lastRow; $r++) {
        $idx = $r - 1
        $f = $origFormula[$c][$idx]
        $v = $origValue[$c][$idx]
        $cell = $ws.Cells.Item($r, $destCol)
        if ($f.Substring(0, 1) -eq "=") {
            $cell.Formula = $f
        } else {
            $cell.Value = $v
        }
    }
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 4: populate the freshly freed column D: header + formula, formatted
# like column A (integer display on the orange "input" fill).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Resize($lastRow, 1).Copy() | Out-Null
$ws.Cells.Item(1, 4).Resize($lastRow, 1).PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 4).Value = "Tick duration (ns)"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=Tableau1[[#This Row],[Tick duration (s)]]*1000000000"
}

# ---------------------------------------------------------------------------
# Step 5: highlight the whole "Prescaler idx = 5" row (worksheet row 6) with
# a yellow fill across the full table width (A:J), marking the chosen
# prescaler configuration.
# ---------------------------------------------------------------------------
$ws.Range("A6:J6").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Step 6: move the active selection to D6, mirroring the author's click on
# the new column while reviewing the highlighted row.
# ---------------------------------------------------------------------------
$ws.Range("D6").Select() | Out-Null
